# Config.xlsx edit: add two new settings rows ("MaxExecutionAttemptsHigh" and
# "RetryIntervalLow") to the "Constants" sheet, right after the existing
# "MaxExecutionAttempts"-family rows, pushing the rows below them down.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Constants")
$ws.Activate()

# Insert a new row at row 7 (pushes the old "MaxLockTimeout" row, previously
# row 7, down to row 8; "RetryInterval" moves from row 8 to row 9, etc.)
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "MaxExecutionAttemptsHigh"
$ws.Range("B7").Value = 99999
$ws.Range("C7").Value = "Maximum number of execution attempts for a process step which by default is high."

# Insert another new row at (current) row 10, right after "RetryInterval"
# (now row 9), pushing the formerly-blank separator row down to row 11.
$ws.Rows.Item(10).Insert()
$ws.Range("A10").Value = "RetryIntervalLow"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Duration, in seconds, between re-execution attempts"

# Restore the cursor/selection state recorded in the saved workbook. The
# "Settings" sheet keeps a remembered selection even though it is not the
# front-most tab, so visit it first and re-activate "Constants" afterwards
# so it remains the active/front tab, matching the saved file.
$ws1 = $wb.Worksheets.Item("Settings")
$ws1.Activate()
$ws1.Range("A24").Select()

$ws.Activate()
$ws.Rows.Item(10).Select()
